$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 5101.5557
$ws.Range("I9").Value = 6252.5293
$ws.Range("J9").Value = 3144.9
$ws.Range("K9").Value = 6252.5293
$ws.Range("L9").Value = 3144.9
$ws.Range("M9").Value = -6083.5293
$ws.Range("N9").Value = -3482.9

$ws.Range("H12").Value = 369.54544
$ws.Range("I12").Value = 438
$ws.Range("K12").Value = 438
$ws.Range("M12").Value = -268

$ws.Range("N19").ClearContents()
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0

$ws.Range("M76").ClearContents()
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0

$ws.Range("M79").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0

$ws.Range("H101").Value = 1866.3334
$ws.Range("I101").Value = 299
$ws.Range("K101").Value = 897
$ws.Range("M101").Value = 725

$ws.Range("H132").Value = 4184.727
$ws.Range("I132").Value = 4099.871
$ws.Range("K132").Value = 12299.613
$ws.Range("M132").Value = -9769.613000000001

$ws.Range("H137").Value = 1663.4286
$ws.Range("I137").Value = 1696.6
$ws.Range("K137").Value = 5089.799999999999
$ws.Range("M137").Value = -2539.799999999999

$ws.Range("H138").Value = 6745.778
$ws.Range("I138").Value = 1227.8334
$ws.Range("K138").Value = 3683.5002
$ws.Range("M138").Value = 1456.4998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 6166.3335
$ws.Range("I28").Value = 6166.3335
$ws.Range("K28").Value = 6166.3335
$ws.Range("M28").Value = -5974.3335

$ws.Range("H45").Value = 3762.125
$ws.Range("I45").Value = 1619.8
$ws.Range("J45").Value = 7332.6665
$ws.Range("K45").Value = 1619.8
$ws.Range("L45").Value = 7332.6665
$ws.Range("M45").Value = -1242.8
$ws.Range("N45").Value = -8086.6665

$ws.Range("M63").ClearContents()
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0

$ws.Range("M66").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0

$ws.Range("H99").Value = 6166.3335
$ws.Range("I99").Value = 6166.3335
$ws.Range("K99").Value = 6166.3335
$ws.Range("M99").Value = -3171.3335

$ws.Range("H132").Value = 25043652
$ws.Range("J132").Value = 71530210
$ws.Range("L132").Value = 214590630
$ws.Range("N132").Value = -214595690

$ws.Range("N133").ClearContents()
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 49999.332
$ws.Range("J95").Value = 49999.332
$ws.Range("L95").Value = 49999.332
$ws.Range("N95").Value = -55491.332

$ws.Range("H134").Value = 2606.5208
$ws.Range("I134").Value = 2823.8684
$ws.Range("K134").Value = 8471.6052
$ws.Range("M134").Value = -5936.6052

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8933639
$ws.Range("I31").Value = 3582
$ws.Range("K31").Value = 3582
$ws.Range("M31").Value = -3287

$ws.Range("H34").Value = 8933639
$ws.Range("I34").Value = 3582
$ws.Range("K34").Value = 3582
$ws.Range("M34").Value = -3380

$ws.Range("H96").Value = 25690.375
$ws.Range("J96").Value = 25690.375
$ws.Range("L96").Value = 25690.375
$ws.Range("N96").Value = -31182.375

$ws.Range("H99").Value = 9069.275
$ws.Range("I99").Value = 6306
$ws.Range("K99").Value = 6306
$ws.Range("M99").Value = -4808

$ws.Range("H126").Value = 9069.275
$ws.Range("I126").Value = 6306
$ws.Range("K126").Value = 18918
$ws.Range("M126").Value = -16448

$ws.Range("H131").Value = 23278.4
$ws.Range("J131").Value = 25848
$ws.Range("L131").Value = 25848
$ws.Range("N131").Value = -35928

$ws.Range("H132").Value = 76616.516
$ws.Range("I132").Value = 97588.52
$ws.Range("J132").Value = 3214.5
$ws.Range("K132").Value = 292765.56
$ws.Range("L132").Value = 9643.5
$ws.Range("M132").Value = -290235.56
$ws.Range("N132").Value = -14703.5

$ws.Range("H141").Value = 100376.6
$ws.Range("J141").Value = 115220.75
$ws.Range("L141").Value = 115220.75
$ws.Range("N141").Value = -125580.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 47686
$ws.Range("I2").Value = 54.5
$ws.Range("J2").Value = 200106.8
$ws.Range("K2").Value = 327
$ws.Range("L2").Value = 1200640.8
$ws.Range("M2").Value = -214
$ws.Range("N2").Value = -1200866.8

$ws.Range("N22").ClearContents()
$ws.Range("H22").Value = 795.6667
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0

$ws.Range("M25").ClearContents()
$ws.Range("H25").Value = 10000
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 10000
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 30000
$ws.Range("N25").Value = -30338

$ws.Range("N27").ClearContents()
$ws.Range("H27").Value = 795.6667
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0

$ws.Range("M30").ClearContents()
$ws.Range("H30").Value = 10000
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 10000
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 30000
$ws.Range("N30").Value = -30204

$ws.Range("H129").Value = 20834488
$ws.Range("I129").Value = 33333984
$ws.Range("J129").Value = 15152899
$ws.Range("K129").Value = 100001952
$ws.Range("L129").Value = 45458697
$ws.Range("M129").Value = -99996952
$ws.Range("N129").Value = -45468697

$ws.Range("H137").Value = 1617
$ws.Range("J137").Value = 3022
$ws.Range("L137").Value = 9066
$ws.Range("N137").Value = -19266

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M95").ClearContents()
$ws.Range("H95").Value = 24999
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 24999
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 24999
$ws.Range("N95").Value = -30491

$ws.Range("H99").Value = 17213.666
$ws.Range("I99").Value = 17213.666
$ws.Range("K99").Value = 17213.666
$ws.Range("M99").Value = -14967.666

$ws.Range("H102").Value = 9999
$ws.Range("I102").Value = 9999
$ws.Range("J102").Value = 9999
$ws.Range("K102").Value = 9999
$ws.Range("L102").Value = 9999
$ws.Range("M102").Value = -8377
$ws.Range("N102").Value = -13243

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3509.125
$ws.Range("I40").Value = 3086.8333
$ws.Range("J40").Value = 4776
$ws.Range("K40").Value = 3086.8333
$ws.Range("L40").Value = 4776
$ws.Range("M40").Value = -2950.8333
$ws.Range("N40").Value = -5048

$ws.Range("H132").Value = 2605.375
$ws.Range("I132").Value = 2112.4
$ws.Range("K132").Value = 6337.200000000001
$ws.Range("M132").Value = -3807.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7110.8887
$ws.Range("J81").Value = 7714.143
$ws.Range("L81").Value = 15428.286
$ws.Range("N81").Value = -17550.286

$ws.Range("H84").Value = 7110.8887
$ws.Range("J84").Value = 7714.143
$ws.Range("L84").Value = 77141.42999999999
$ws.Range("N84").Value = -87749.42999999999

$ws.Range("H100").Value = 101002280
$ws.Range("I100").Value = 144288420
$ws.Range("K100").Value = 288576840
$ws.Range("M100").Value = -288576299

$ws.Range("H113").Value = 1983.8
$ws.Range("I113").Value = 1983.8
$ws.Range("K113").Value = 5951.4
$ws.Range("M113").Value = -3781.4

$ws.Range("H133").Value = 71514.8
$ws.Range("J133").Value = 71514.8
$ws.Range("L133").Value = 71514.8
$ws.Range("N133").Value = -81634.8

$ws.Range("H136").Value = 3446.2778
$ws.Range("I136").Value = 2766.647
$ws.Range("J136").Value = 15000
$ws.Range("K136").Value = 8299.940999999999
$ws.Range("L136").Value = 45000
$ws.Range("M136").Value = -5749.940999999999
$ws.Range("N136").Value = -50100
